$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value (45188 = 2023-09-19).
# For every data row (2 through 230) bump that date forward by one day
# to 45189 (2023-09-20).
for ($row = 2; $row -le 230; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45188) {
        $cell.Value2 = 45189
    }
}
